$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.045.68"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").Value = "2.268.80"
$ws.Range("E3").Value = "  -4.25%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'488.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.27%  "

$ws.Range("D6").Value = "'127.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.94%  "

$ws.Range("D9").Value = "2.273.93"
$ws.Range("E9").Value = "  -4.22%  "

$ws.Range("D10").Value = "'0.0926"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.59%  "

$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Value = "'4.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.66%  "

$ws.Range("E13").Value = "  -3.22%  "

$ws.Range("D14").Value = "2.676.00"
$ws.Range("E14").Value = "  -4.00%  "

$ws.Range("D15").Value = "'21.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").Value = "54.030.72"
$ws.Range("E16").Value = "  -3.70%  "

$ws.Range("E17").Value = "  -2.63%  "

$ws.Range("D18").Value = "2.287.08"
$ws.Range("E18").Value = "  -5.38%  "

$ws.Range("D19").Value = "'3.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  -4.08%  "

$ws.Range("D21").Value = "'302.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("E22").Value = "  -1.81%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'63.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").Value = "'0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'0.366"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").Value = "'7.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("D29").Value = "'169.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("D31").Value = "'1.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").Value = "'0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").Value = "'5.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "'17.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  -1.25%  "

$ws.Range("D38").Value = "'0.836"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.22%  "

$ws.Range("D39").Value = "'3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "'35.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.90%  "

$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "'122.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.18%  "

$ws.Range("D45").Value = "'4.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("D46").Value = "'0.0879"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.74%  "

$ws.Range("E47").Value = "  -4.16%  "

$ws.Range("D48").Value = "'239.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("D50").Value = "'0.0203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("D51").Value = "'16.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.81%  "
